$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 35
$prev = $row - 1

# Copy the formatting (styles) of the previous row so the new row matches
# the existing look (bold/bordered index column, date-formatted date column).
$ws.Range("A$prev`:Y$prev").Copy()
$ws.Range("A$row").PasteSpecial(-4122)

# Yesterday's data appended as a new row.
$values = @{
    A = 33
    B = 43924
    C = 1095917
    D = 225796
    E = 58787
    F = 275586
    G = 9707
    H = 7087
    I = 102987
    J = 0
    K = 2935
    L = 3946
    M = 1287
    N = 139
    O = 164
    P = 93
    Q = 1363
    R = 1620
    S = 0
    T = 0
    U = 0
    V = 2
    W = 1
    X = 0
    Y = 51
}

foreach ($col in $values.Keys) {
    $ws.Range("$col$row").Value = $values[$col]
}
